$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.014.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.80%  "

$ws.Range("D3").Value = "'1.653.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.02%  "

$ws.Range("E4").Value = "  -0.37%  "

$ws.Range("D5").Value = "'216.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.87%  "

$ws.Range("D6").Value = "'0.5211"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("E7").Value = "  -0.38%  "

$ws.Range("D8").Value = "'0.2621"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.95%  "

$ws.Range("D9").Value = "'0.06269"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.53%  "

$ws.Range("D10").Value = "'20.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.20%  "

$ws.Range("D11").Value = "'0.07730"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.11%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.472"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.18%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.647.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.61%  "

$ws.Range("D14").Value = "'1.879.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.07%  "

$ws.Range("D15").Value = "'0.5432"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.28%  "

$ws.Range("D16").Value = "'0.0₅8094"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.42%  "

$ws.Range("D17").Value = "'64.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("D18").Value = "'26.018.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.86%  "

$ws.Range("D19").Value = "'1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.34%  "

$ws.Range("D20").Value = "'4.578"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.12%  "

$ws.Range("D21").Value = "'191.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.21%  "

$ws.Range("D22").Value = "'10.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.55%  "

$ws.Range("D23").Value = "'5.987"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.31%  "

$ws.Range("E24").Value = "  -0.46%  "

$ws.Range("D25").Value = "'138.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.41%  "

$ws.Range("D26").Value = "'0.1233"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.86%  "

$ws.Range("D27").Value = "'7.247"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.51%  "

$ws.Range("D28").Value = "'16.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.29%  "

$ws.Range("D29").Value = "'1.399"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.13%  "

$ws.Range("D30").Value = "'0.05952"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.81%  "

$ws.Range("D31").Value = "'1.273"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.92%  "

$ws.Range("D32").Value = "'3.507"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.18%  "

$ws.Range("D33").Value = "'3.247"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.20%  "

$ws.Range("D34").Value = "'1.561"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.49%  "

$ws.Range("D35").Value = "'0.9507"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.43%  "

$ws.Range("D36").Value = "'2.411"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("D37").Value = "'2.751"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.87%  "

$ws.Range("D38").Value = "'0.5690"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.46%  "

$ws.Range("D39").Value = "'0.01595"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.09%  "

$ws.Range("D40").Value = "'5.902"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.95%  "

$ws.Range("D41").Value = "'0.8469"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.74%  "

$ws.Range("E42").Value = "  -0.16%  "

$ws.Range("D43").Value = "'100.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.75%  "

$ws.Range("D44").Value = "'1.002.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.59%  "

$ws.Range("D45").Value = "'1.794.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.17%  "

$ws.Range("E46").Value = "  -1.00%  "

$ws.Range("D47").Value = "'56.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.96%  "

$ws.Range("D48").Value = "'0.9995"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.35%  "

$ws.Range("D49").Value = "'7.959"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.25%  "

$ws.Range("D50").Value = "'0.4296"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.55%  "

$ws.Range("D51").Value = "'1.477"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.60%  "
